$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Columns used in the sheet: A .. AY (1 .. 51)
$lastCol = 51
$stagingRowBase = 500

function Get-RowRange([int]$r) {
    return $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, $lastCol))
}

function Swap-Rows([int]$r1, [int]$r2, [int]$stagingRow) {
    $rng1 = Get-RowRange $r1
    $rng2 = Get-RowRange $r2
    $stage = Get-RowRange $stagingRow

    # stage = row r1 (save a copy of the original row r1 contents)
    $rng1.Copy() | Out-Null
    $stage.PasteSpecial(-4104) | Out-Null   # xlPasteAll

    # row r1 = row r2
    $rng2.Copy() | Out-Null
    $rng1.PasteSpecial(-4104) | Out-Null

    # row r2 = stage (original row r1)
    $stage.Copy() | Out-Null
    $rng2.PasteSpecial(-4104) | Out-Null

    # clear the staging area
    $stage.ClearContents() | Out-Null
    $excel.CutCopyMode = $false
}

# Rows 43 and 44 swap their entire contents.
# Row 43 uniquely had AJ/AK/AO populated; after the swap those move to row 44,
# so row 43 must end up without them.
Swap-Rows 43 44 $stagingRowBase
$ws.Cells.Item(44, 36).ClearContents() | Out-Null  # AJ44
$ws.Cells.Item(44, 37).ClearContents() | Out-Null  # AK44
$ws.Cells.Item(44, 41).ClearContents() | Out-Null  # AO44

# Rows 55 and 56 swap their entire contents (no extra columns to fix up).
Swap-Rows 55 56 $stagingRowBase

# Rows 66 and 67 swap their entire contents.
# Row 67 uniquely had I/J populated ("75" / "bålar"); after the swap those
# move to row 66, so row 67 must end up without them.
Swap-Rows 66 67 $stagingRowBase
$ws.Cells.Item(67, 9).ClearContents() | Out-Null   # I67
$ws.Cells.Item(67, 10).ClearContents() | Out-Null  # J67

# Rows 76, 77, 78 undergo a 3-way rotation of their content:
#   new row 76 = old row 77
#   new row 77 = old row 78
#   new row 78 = old row 76
$r76 = Get-RowRange 76
$r77 = Get-RowRange 77
$r78 = Get-RowRange 78
$stage76 = Get-RowRange ($stagingRowBase + 3)

# Save old row 76 into staging
$r76.Copy() | Out-Null
$stage76.PasteSpecial(-4104) | Out-Null

# row76 = old row77
$r77.Copy() | Out-Null
$r76.PasteSpecial(-4104) | Out-Null

# row77 = old row78
$r78.Copy() | Out-Null
$r77.PasteSpecial(-4104) | Out-Null

# row78 = old row76 (from staging)
$stage76.Copy() | Out-Null
$r78.PasteSpecial(-4104) | Out-Null

$stage76.ClearContents() | Out-Null
$excel.CutCopyMode = $false
